$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Append the 4 new daily rows (2025-12-07 .. 2025-12-10) after the
# existing last row (62, date 2025-12-06).
$newRows = @(
    @{ Row = 63; Date = "2025-12-07"; NonHttps = 0.0; Https = 25.0 },
    @{ Row = 64; Date = "2025-12-08"; NonHttps = 0.0; Https = 26.0 },
    @{ Row = 65; Date = "2025-12-09"; NonHttps = 0.0; Https = 27.0 },
    @{ Row = 66; Date = "2025-12-10"; NonHttps = 0.0; Https = 27.0 }
)

foreach ($r in $newRows) {
    $dateCell = $ws.Cells.Item($r.Row, 1)
    # Force the date-like text to be stored as a literal text string
    # (matching the rest of the column) instead of being auto-converted
    # to a date serial number.
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $r.Date
    $dateCell.ClearFormats()

    $ws.Cells.Item($r.Row, 2).Value = $r.NonHttps
    $ws.Cells.Item($r.Row, 3).Value = $r.Https
}
